# Update "想去人数" (column F) counts on the "展览", "演出" and "全部类型"
# sheets to the freshly scraped values (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1030
$ws1.Range("F3").Value = 13545
$ws1.Range("F5").Value = 1029
$ws1.Range("F6").Value = 20
$ws1.Range("F7").Value = 1740
$ws1.Range("F8").Value = 141
$ws1.Range("F12").Value = 32
$ws1.Range("F13").Value = 13556
$ws1.Range("F15").Value = 602
$ws1.Range("F16").Value = 8964
$ws1.Range("F18").Value = 8056
$ws1.Range("F19").Value = 255
$ws1.Range("F20").Value = 12
$ws1.Range("F21").Value = 149
$ws1.Range("F28").Value = 19
$ws1.Range("F30").Value = 208
$ws1.Range("F31").Value = 184

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 40

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1030
$ws4.Range("F3").Value = 13545
$ws4.Range("F5").Value = 1029
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 1740
$ws4.Range("F8").Value = 141
$ws4.Range("F12").Value = 32
$ws4.Range("F13").Value = 13556
$ws4.Range("F15").Value = 602
$ws4.Range("F16").Value = 8964
$ws4.Range("F18").Value = 8056
$ws4.Range("F19").Value = 255
$ws4.Range("F20").Value = 12
$ws4.Range("F21").Value = 149
$ws4.Range("F28").Value = 19
$ws4.Range("F29").Value = 40
$ws4.Range("F32").Value = 208
$ws4.Range("F33").Value = 184
